$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "-"
$ws.Range("C2").Value = "-"
$ws.Range("B4").Value = "-"
$ws.Range("D14").Value = "-"
$ws.Range("F15").Value = "-"
$ws.Range("B18").Value = "[Aderci-Fresagem-2NA, Leonardo-M.Maq.E.I.-2NA, Ismail-Metrologia 2-2NA, Aderci-CAD/CAM-2NA]"
$ws.Range("C18").Value = "[João Bosco-Fundição-2NA, Leandro-M.S.R.A.C.-2NA, Paulo Rob.-Usin. CNC-2NA, Euclides-Soldagem-2NA]"
$ws.Range("D18").Value = "[Paulo Rob.-Usin. CNC-2NA, Leandro-M.S.R.A.C.-2NA, João Bosco-Fundição-2NA, Euclides-Soldagem-2NA]"
$ws.Range("E18").Value = "[Guilherme-C.L.P.-2NA, Leonardo-Retífica-2NA, Guilherme-C.pneumática-2NA, Guilherme-C. Hidráulica-2NA]"
$ws.Range("F18").Value = "-"
$ws.Range("B19").Value = "[Aderci-Fresagem-2NA, Leonardo-M.Maq.E.I.-2NA, Ismail-Metrologia 2-2NA, Aderci-CAD/CAM-2NA]"
$ws.Range("C19").Value = "[João Bosco-Fundição-2NA, Leandro-M.S.R.A.C.-2NA, Paulo Rob.-Usin. CNC-2NA, Euclides-Soldagem-2NA]"
$ws.Range("D19").Value = "[Paulo Rob.-Usin. CNC-2NA, Leandro-M.S.R.A.C.-2NA, João Bosco-Fundição-2NA, Euclides-Soldagem-2NA]"
$ws.Range("E19").Value = "[Guilherme-C.L.P.-2NA, Leonardo-Retífica-2NA, Guilherme-C.pneumática-2NA, Guilherme-C. Hidráulica-2NA]"
$ws.Range("B20").Value = "[Aderci-Fresagem-2NA, Leonardo-M.Maq.E.I.-2NA, Ismail-Metrologia 2-2NA, Aderci-CAD/CAM-2NA]"
$ws.Range("C20").Value = "Claudinei-Des. Maq. Cad-"
$ws.Range("D20").Value = "Euclides-Gest. Int.-"
$ws.Range("E20").Value = "[Guilherme-C.L.P.-2NA, Leonardo-Retífica-2NA, Guilherme-C.pneumática-2NA, Guilherme-C. Hidráulica-2NA]"
$ws.Range("F20").Value = "Claudinei-Des. Maq. Cad-"
$ws.Range("B21").Value = "[Aderci-Fresagem-2NA, Leonardo-M.Maq.E.I.-2NA, Ismail-Metrologia 2-2NA, Aderci-CAD/CAM-2NA]"
$ws.Range("C21").Value = "Claudinei-Des. Maq. Cad-"
$ws.Range("D21").Value = "Euclides-Gest. Int.-"
$ws.Range("E21").Value = "[Guilherme-C.L.P.-2NA, Leonardo-Retífica-2NA, Guilherme-C.pneumática-2NA, Guilherme-C. Hidráulica-2NA]"
$ws.Range("F21").Value = "Claudinei-Elemaq.-"
